# "add counter for existing file"
#
# The existing "Serials"/"Movie" rules pointed at local drive paths
# (F:\Serials\_Unknow, F:\Filmy) that may not exist/be mounted on every
# machine running the DownloadFolderManager worker. Re-point the "already
# exists" target folders at the network share (\\denynaspc\...) instead,
# and likewise normalise the two *.ext style patterns (doc-vyplatni paska /
# doc-vodnestocne rows) to plain ".pdf"/".docx" extensions plus network
# share destinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("migrate_rule")

# Row 4 ("Serials"): target folder -> network share, with a real hyperlink
# (Excel auto-applies the built-in Hyperlink cell style when this happens).
$ws.Range("E4").Value = "\\denynaspc\Serials\_Unknow"
$ws.Hyperlinks.Add($ws.Range("E4"), "\\denynaspc\Serials\_Unknow")

# Row 6 ("Movie"): target folder -> network share, with a real hyperlink.
$ws.Range("E6").Value = "\\denynaspc\Filmy\Filmy"
$ws.Hyperlinks.Add($ws.Range("E6"), "\\denynaspc\Filmy\Filmy")

# Row 7 ("Movie subtitle"): same target folder as row 6, but plain text.
$ws.Range("E7").Value = "\\denynaspc\Filmy\Filmy"

# Row 14 ("doc-vyplatni paska"): extension pattern + target folder.
$ws.Range("E14").Value = "\\denynaspc\Me\Document\Payslip"
$ws.Range("C14").Value = ".pdf"

# Row 15 ("doc-vodnestocne"): extension pattern + target folder.
$ws.Range("C15").Value = ".docx"
$ws.Range("E15").Value = "\\denynaspc\Me\Document\Documents\Invoice\Paid\Watter"

# Leave the cursor where the author left it when they saved.
[void]$ws.Range("I26").Select()

Write-Output "migrate_rule rules updated"
